# Atualizado por script em 02-11-2023 14:45
#
# 1) Rows 20 and 21 (match data, columns F:V) were swapped: the
#    "Al Bataeh vs Shabab Al-Ahli Dubai" match moves from row 20 to row 21,
#    and the "Al Wahda vs Hatta" match moves from row 21 to row 20.
#    Columns A:E (Indice / pais / torneio / temporada / data_partida) stay put.
# 2) A new match row (44) is appended: Al Jazira vs Al Sharjah.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20 becomes the old row 21 data (Al Wahda vs Hatta) ---
$ws.Range("F20").Value = "Al Wahda"
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = "Hatta"
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 1.36
$ws.Range("K20").Value = "19/09/2023 16:42"
$ws.Range("L20").Value = 1.24
$ws.Range("M20").Value = "24/09/2023 14:22"
$ws.Range("N20").Value = 5.48
$ws.Range("O20").Value = "19/09/2023 16:42"
$ws.Range("P20").Value = 6.54
$ws.Range("Q20").Value = "24/09/2023 15:15"
$ws.Range("R20").Value = 6.99
$ws.Range("S20").Value = "19/09/2023 16:42"
$ws.Range("T20").Value = 10.02
$ws.Range("U20").Value = "24/09/2023 15:15"
$ws.Range("V20").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-wahda-hatta/pWp8qWSh/"

# --- Row 21 becomes the old row 20 data (Al Bataeh vs Shabab Al-Ahli Dubai) ---
$ws.Range("F21").Value = "Al Bataeh"
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = "Shabab Al-Ahli Dubai"
$ws.Range("I21").Value = 2
$ws.Range("J21").Value = 6.4
$ws.Range("K21").Value = "17/09/2023 15:42"
$ws.Range("L21").Value = 6.67
$ws.Range("M21").Value = "24/09/2023 15:15"
$ws.Range("N21").Value = 5.03
$ws.Range("O21").Value = "17/09/2023 15:42"
$ws.Range("P21").Value = 4.85
$ws.Range("Q21").Value = "24/09/2023 15:15"
$ws.Range("R21").Value = 1.37
$ws.Range("S21").Value = "17/09/2023 15:42"
$ws.Range("T21").Value = 1.44
$ws.Range("U21").Value = "24/09/2023 15:15"
$ws.Range("V21").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-bataeh-shabab-al-ahli-dubai/6un4pCDn/"

# --- Append new row 44 (Al Jazira vs Al Sharjah) ---
# Copy formatting from the last existing row (43) first so the new row
# picks up the same per-column styles (bold/bordered index in A, the
# custom datetime format in E, plain for everything else).
$ws.Range("A43:V43").Copy()
$ws.Range("A44:V44").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = "united-arab-emirates"
$ws.Range("C44").Value = "uae-league"
$ws.Range("D44").Value = "2023-2024"
$ws.Range("E44").Value = 45232.57291666666
$ws.Range("F44").Value = "Al Jazira"
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = "Al Sharjah"
$ws.Range("I44").Value = 1
$ws.Range("J44").Value = 2.05
$ws.Range("K44").Value = "30/10/2023 18:42"
$ws.Range("L44").Value = 2.45
$ws.Range("M44").Value = "02/11/2023 13:43"
$ws.Range("N44").Value = 3.8
$ws.Range("O44").Value = "30/10/2023 18:42"
$ws.Range("P44").Value = 3.85
$ws.Range("Q44").Value = "02/11/2023 13:43"
$ws.Range("R44").Value = 3.3
$ws.Range("S44").Value = "30/10/2023 18:42"
$ws.Range("T44").Value = 2.65
$ws.Range("U44").Value = "02/11/2023 13:43"
$ws.Range("V44").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-jazira-al-sharjah/tt7G0O3p/"

Write-Output "edit complete"
